$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update country labels whose rank order changed ---
# Iran / Brasil swapped position (Brasil overtook Iran in case count)
$ws.Range("A12").Value = "Brasil"
$ws.Range("A13").Value = "Iran"

# Uruguay / Burkina Faso swapped position
$ws.Range("A105").Value = "Uruguay"
$ws.Range("A106").Value = "Burkina Faso"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1157421
$ws.Range("C4").Value = 26391
$ws.Range("D4").Value = 160552
$ws.Range("E4").Value = 929621
$ws.Range("F4").Value = 16455
$ws.Range("G4").Value = 1495
$ws.Range("H4").Value = 67248

# --- Row 12: now Brasil ---
$ws.Range("B12").Value = 96559
$ws.Range("C12").Value = 4450
$ws.Range("D12").Value = 40937
$ws.Range("E12").Value = 48872
$ws.Range("F12").Value = 8318
$ws.Range("G12").Value = 340
$ws.Range("H12").Value = 6750

# --- Row 13: now Iran ---
$ws.Range("B13").Value = 96448
$ws.Range("C13").Value = 802
$ws.Range("D13").Value = 77350
$ws.Range("E13").Value = 12942
$ws.Range("F13").Value = 2787
$ws.Range("G13").Value = 65
$ws.Range("H13").Value = 6156

# --- Row 48 ---
$ws.Range("B48").Value = 7285
$ws.Range("C48").Value = 279
$ws.Range("D48").Value = 1666
$ws.Range("E48").Value = 5295
$ws.Range("F48").Value = 118
$ws.Range("G48").Value = 10
$ws.Range("H48").Value = 324

# --- Row 105: now Uruguay ---
$ws.Range("B105").Value = 652
$ws.Range("C105").Value = 4
$ws.Range("D105").Value = 440
$ws.Range("E105").Value = 195
$ws.Range("F105").Value = 10
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 17

# --- Row 106: now Burkina Faso ---
$ws.Range("B106").Value = 652
$ws.Range("C106").Value = 3
$ws.Range("D106").Value = 535
$ws.Range("E106").Value = 73
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 44

# --- Row 152 ---
$ws.Range("B152").Value = 114
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 51
$ws.Range("E152").Value = 56
$ws.Range("F152").Value = 4
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 7
